# This workbook's only worksheet is protected, so it must be unprotected
# before any cell values can be changed, and re-protected afterwards to
# keep the workbook's overall protected state intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A11):
# 2021-03-18 -> 2021-03-19
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-19 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) figures for each holding row.
$ws.Range("D2").Value = 0.4988340920510653
$ws.Range("E2").Value = -0.005095541401273884

$ws.Range("D3").Value = 0.2394599414087377
$ws.Range("E3").Value = 0.002834645669291369

$ws.Range("D4").Value = 0.09801208933713243
$ws.Range("E4").Value = 0.005028365136668489

$ws.Range("D5").Value = 0.1037432414948896
$ws.Range("E5").Value = -0.004156999226604841

$ws.Range("D6").Value = 0.03127341218381328
$ws.Range("E6").Value = 0.0004810467577449629

$ws.Range("D7").Value = 0.02867722352436171
$ws.Range("E7").Value = 0.004616449949988333

$ws.Range("D8").Value = 1
$ws.Range("E8").Value = -0.001654034743555521

# Restore sheet protection to match the workbook's original protected state.
$ws.Protect()
